$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.839.56"
$ws.Range("E2").Value = "  +11.31%  "

$ws.Range("D3").Value = "1.745.69"
$ws.Range("E3").Value = "  +7.29%  "

$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  +0.90%  "

$ws.Range("D5").Value = "315.87"
$ws.Range("E5").Value = "  +4.17%  "

$ws.Range("D6").Value = "0.9937"
$ws.Range("E6").Value = "  +1.06%  "

$ws.Range("D7").Value = "0.3827"
$ws.Range("E7").Value = "  +4.06%  "

$ws.Range("D8").Value = "0.3638"
$ws.Range("E8").Value = "  +5.99%  "

$ws.Range("D9").Value = "50.70"
$ws.Range("E9").Value = "  +18.80%  "

$ws.Range("D10").Value = "1.227"
$ws.Range("E10").Value = "  +6.06%  "

$ws.Range("D11").Value = "0.07707"
$ws.Range("E11").Value = "  +8.46%  "

$ws.Range("D12").Value = "0.9951"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("D13").Value = "21.78"
$ws.Range("E13").Value = "  +7.40%  "

$ws.Range("D14").Value = "6.455"
$ws.Range("E14").Value = "  +8.93%  "

$ws.Range("D15").Value = "7.068"
$ws.Range("E15").Value = "  +6.20%  "

$ws.Range("D16").Value = "1.741.19"
$ws.Range("E16").Value = "  +7.44%  "

$ws.Range("D17").Value = "0.00001157"
$ws.Range("E17").Value = "  +6.84%  "

$ws.Range("D18").Value = "0.9935"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").Value = "0.06827"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("D20").Value = "87.11"
$ws.Range("E20").Value = "  +10.97%  "

$ws.Range("E21").Value = "  +8.30%  "

$ws.Range("D22").Value = "6.468"
$ws.Range("E22").Value = "  +7.26%  "

$ws.Range("D23").Value = "12.76"
$ws.Range("E23").Value = "  +8.52%  "

$ws.Range("D24").Value = "25.763.35"
$ws.Range("E24").Value = "  +11.28%  "

$ws.Range("D25").Value = "2.432"
$ws.Range("E25").Value = "  +2.49%  "

$ws.Range("D26").Value = "2.938"
$ws.Range("E26").Value = "  +12.32%  "

$ws.Range("D27").Value = "20.70"
$ws.Range("E27").Value = "  +6.46%  "

$ws.Range("D28").Value = "154.21"
$ws.Range("E28").Value = "  +2.48%  "

$ws.Range("D29").Value = "134.32"
$ws.Range("E29").Value = "  +7.27%  "

$ws.Range("D30").Value = "1.933.55"
$ws.Range("E30").Value = "  +7.59%  "

$ws.Range("E31").Value = "  +21.95%  "

$ws.Range("D32").Value = "7.044"
$ws.Range("E32").Value = "  +15.49%  "

$ws.Range("D33").Value = "4.344"
$ws.Range("E33").Value = "  +6.90%  "

$ws.Range("D34").Value = "14.32"
$ws.Range("E34").Value = "  +19.51%  "

$ws.Range("D35").Value = "1.804"
$ws.Range("E35").Value = "  +7.87%  "

$ws.Range("D36").Value = "0.08695"
$ws.Range("E36").Value = "  +5.11%  "

$ws.Range("D37").Value = "5.644"
$ws.Range("E37").Value = "  +8.10%  "

$ws.Range("D38").Value = "0.06719"
$ws.Range("E38").Value = "  +8.27%  "

$ws.Range("D39").Value = "9.301"

$ws.Range("E40").Value = "  +9.80%  "

$ws.Range("D41").Value = "0.2221"
$ws.Range("E41").Value = "  +9.51%  "

$ws.Range("D42").Value = "1.299"
$ws.Range("E42").Value = "  +3.31%  "

$ws.Range("E43").Value = "  +9.97%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "0.9932"
$ws.Range("E44").Value = "  +1.09%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  +6.76%  "

$ws.Range("D46").Value = "0.6347"
$ws.Range("E46").Value = "  +9.52%  "

$ws.Range("D47").Value = "3.899"
$ws.Range("E47").Value = "  +2.53%  "

$ws.Range("D48").Value = "2.177"
$ws.Range("E48").Value = "  +9.20%  "

$ws.Range("D49").Value = "131.49"
$ws.Range("E49").Value = "  +3.77%  "

$ws.Range("D50").Value = "0.07486"
$ws.Range("E50").Value = "  +7.61%  "

$ws.Range("D51").Value = "79.39"
$ws.Range("E51").Value = "  +6.70%  "
